$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-92). All of them move from serial 45177 (2023-09-08) to
# serial 45178 (2023-09-09).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row()

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}
